$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking Price (column D) values stay as literal text,
# matching the original inline-string cell type (avoid Excel auto-converting
# "581.70" style values into floating point numbers like 581.70000000000005).
$priceCells = @(
    "D2",
    "D3",
    "D5",
    "D6",
    "D7",
    "D10",
    "D11",
    "D12",
    "D13",
    "D14",
    "D15",
    "D16",
    "D18",
    "D19",
    "D21",
    "D22",
    "D23",
    "D24",
    "D25",
    "D26",
    "D28",
    "D30",
    "D32",
    "D33",
    "D34",
    "D36",
    "D37",
    "D39",
    "D40",
    "D41",
    "D42",
    "D44",
    "D45",
    "D46",
    "D47",
    "D50",
    "D51"
)
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '66.432.35'
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").Value = '3.531.18'
$ws.Range("E3").Value = '  +0.48%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '581.70'
$ws.Range("E5").Value = '  +5.98%  '
$ws.Range("D6").Value = '179.91'
$ws.Range("E6").Value = '  -6.24%  '
$ws.Range("D7").Value = '0.632'
$ws.Range("E7").Value = '  +4.42%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("E9").Value = '  +1.58%  '
$ws.Range("D10").Value = '0.163'
$ws.Range("E10").Value = '  +8.10%  '
$ws.Range("D11").Value = '55.86'
$ws.Range("E11").Value = '  +1.87%  '
$ws.Range("D12").Value = '0.0000283'
$ws.Range("E12").Value = '  +5.65%  '
$ws.Range("D13").Value = '9.32'
$ws.Range("E13").Value = '  -0.02%  '
$ws.Range("D14").Value = '4.093.20'
$ws.Range("E14").Value = '  +0.70%  '
$ws.Range("D15").Value = '3.536.53'
$ws.Range("E15").Value = '  +0.85%  '
$ws.Range("D16").Value = '18.50'
$ws.Range("E16").Value = '  +1.96%  '
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("D18").Value = '66.400.94'
$ws.Range("E18").Value = '  -0.89%  '
$ws.Range("D19").Value = '12.09'
$ws.Range("E19").Value = '  +1.39%  '
$ws.Range("E20").Value = '  +1.81%  '
$ws.Range("D21").Value = '416.27'
$ws.Range("E21").Value = '  -1.96%  '
$ws.Range("D22").Value = '4.27'
$ws.Range("E22").Value = '  +9.67%  '
$ws.Range("D23").Value = '4.36'
$ws.Range("E23").Value = '  +4.78%  '
$ws.Range("D24").Value = '85.83'
$ws.Range("E24").Value = '  +1.48%  '
$ws.Range("D25").Value = '13.34'
$ws.Range("E25").Value = '  +11.79%  '
$ws.Range("D26").Value = '11.27'
$ws.Range("E26").Value = '  +1.24%  '
$ws.Range("E27").Value = '  -0.56%  '
$ws.Range("D28").Value = '6.05'
$ws.Range("E28").Value = '  -1.66%  '
$ws.Range("E29").Value = '  +3.46%  '
$ws.Range("D30").Value = '30.58'
$ws.Range("E30").Value = '  +1.57%  '
$ws.Range("E31").Value = '  -0.28%  '
$ws.Range("D32").Value = '604.69'
$ws.Range("E32").Value = '  -7.70%  '
$ws.Range("D33").Value = '11.78'
$ws.Range("E33").Value = '  +0.78%  '
$ws.Range("D34").Value = '0.111'
$ws.Range("E34").Value = '  +1.25%  '
$ws.Range("D36").Value = '0.155'
$ws.Range("E36").Value = '  +10.57%  '
$ws.Range("D37").Value = '0.0₃0815'
$ws.Range("E37").Value = '  -0.04%  '
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").Value = '37.43'
$ws.Range("E39").Value = '  -2.71%  '
$ws.Range("D40").Value = '3.63'
$ws.Range("E40").Value = '  +10.31%  '
$ws.Range("D41").Value = '0.387'
$ws.Range("E41").Value = '  -0.13%  '
$ws.Range("D42").Value = '3.254.15'
$ws.Range("E42").Value = '  +8.89%  '
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").Value = '2.95'
$ws.Range("E44").Value = '  +3.21%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = '3.36'
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").Value = '2.58'
$ws.Range("E46").Value = '  -2.14%  '
$ws.Range("D47").Value = '0.0424'
$ws.Range("E47").Value = '  +2.04%  '
$ws.Range("E48").Value = '  +0.01%  '
$ws.Range("E49").Value = '  +2.19%  '
$ws.Range("D50").Value = '8.68'
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("D51").Value = '139.87'
$ws.Range("E51").Value = '  -0.13%  '
